$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7170
$ws.Range("C3").Value = 153067
$ws.Range("C4").Value = 144613
$ws.Range("C5").Value = 8454
$ws.Range("C8").Value = 63.76
